# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
#
# The workbook's worker/period table (rows 16-18) gets a new row inserted
# (one extra "Valor Mora" period for VICTOR BELLIDO RIVERA), the period
# count goes from 2 to 3, the total "Valor Mora" amount is recalculated,
# and a couple of header/label cells are refreshed as part of the same
# resave.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a brand-new row at position 19 (this pushes the old rows 23/24
#    "firma" block down to 24/25, exactly like Excel's own Rows().Insert()).
$ws.Rows(19).Insert()

# 2) Clone row 18 (still holding the old "last data row" content/format)
#    down into the freshly inserted row 19 so the bottom border / shading
#    of the table moves to the new last row.
$ws.Range("B18:J18").Copy()
$ws.Range("B19:J19").PasteSpecial(-4104)

# 3) Row 18 is no longer the last row of the table, so restyle it to match
#    the "middle" rows (copy row 17's formatting only).
$ws.Range("B17:J17").Copy()
$ws.Range("B18:J18").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 4) Refresh the worker / period detail rows.
#    Row 16: CC / 73008904 / LUIS MANUEL HERNANDEZ PASTRANA / periodo 2505
$ws.Range("C16").Value = "73008904"
$ws.Range("D16").Value = "LUIS MANUEL HERNANDEZ PASTRANA"
$ws.Range("E16").Value = "2505"

#    Row 17: CC / 73000468 / VICTOR BELLIDO RIVERA / periodo 2507
$ws.Range("C17").Value = "73000468"
$ws.Range("D17").Value = "VICTOR BELLIDO RIVERA"
$ws.Range("E17").Value = "2507"

#    Row 18: CC / 73000468 / VICTOR BELLIDO RIVERA / periodo 2506
$ws.Range("C18").Value = "73000468"
$ws.Range("D18").Value = "VICTOR BELLIDO RIVERA"
$ws.Range("E18").Value = "2506"

#    Row 19 (new): CC / 73000468 / VICTOR BELLIDO RIVERA / periodo 2505
$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "73000468"
$ws.Range("D19").Value = "VICTOR BELLIDO RIVERA"
$ws.Range("E19").Value = "2505"
$ws.Range("F19").Value = 56940
$ws.Range("G19").Value = 1423500

# 5) Header / summary cells.
$ws.Range("E11").Value = 227760   # VALOR MORA total
$ws.Range("F13").Value = 3        # Cant. Periodos (was 2)
